$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. The old (stray) "_GoBack" bookmark sitting on the very first
#    paragraph gets dropped - Word re-homes it wherever the cursor
#    last was, which after this edit is at the end of the newly
#    typed bullet (handled in step 3 below).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Add one more item to the "What you plan to do next?" list:
#    "We need to find a proper way to cleanup 'NA's in the data"
#    This goes right after the "...improve our overall predictions"
#    bullet, at the same outline level as its sibling bullets
#    (ilvl 0 of numId 8).
# ------------------------------------------------------------------
$lastBullet = $d.Paragraphs.Item(27).Range
$lastBullet.Collapse(0)
$lastBullet.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item(28)
$newPara.Range.ListFormat.ListLevelNumber = 1

$textRange = $newPara.Range
$textRange.InsertAfter("We nee")

$textRange = $d.Paragraphs.Item(28).Range
$textRange.InsertAfter("d to find a proper way to cleanup " + [char]0x2018 + "NA" + [char]0x2019 + "s in the dataZZZ")

# Exclude the paragraph mark so Start/End below line up with the
# visible text only.
[void]$textRange.MoveEnd(1, -1)

# ------------------------------------------------------------------
# 3. Word leaves "_GoBack" collapsed at the spot of the last edit,
#    i.e. right after "...in the data". Bookmarking a genuinely
#    collapsed range directly at end-of-document is unreliable, so
#    bookmark the trailing "ZZZ" placeholder (a real, non-empty
#    range) and then delete that placeholder text - the bookmark
#    collapses in place, exactly where we want it.
# ------------------------------------------------------------------
$placeholder = $d.Range($textRange.End - 3, $textRange.End)
$d.Bookmarks.Add("_GoBack", $placeholder)
$placeholder = $d.Range($textRange.End - 3, $textRange.End)
$placeholder.Text = ""

"done"
